$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: Correspond Handoff Datetime (E4) and Correspond Handback DateTime (H4)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-20 16:58:46"
$wsZh.Range("H4").Value = "2016-03-20 16:59:33"

# "de-de" sheet: Correspond Handoff Datetime (E4) and Correspond Handback DateTime (H4)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-20 16:58:55"
$wsDe.Range("H4").Value = "2016-03-20 16:59:48"
